# Applies the FY2022 "dropped COVID revs" data refresh:
#  - exp_long: 3 corrected per-row expenditure dollar figures
#  - Table 1: re-sorted top rows (Revenue dropped below Toll Highway/Debt
#    Service as its FY22 $ billions fell from 2.2 to 1.8) + refreshed
#    Commerce 1yr/5yr figures
#  - Table 4.b: refreshed Commerce (row 8) and Revenue (row 27) trend figures
#  - year_totals: rebuilt from the recreated allexpfiles22.csv - rows now
#    run most-recent-year-first, and revenue/expenditure totals for several
#    years (notably FY22 & FY21) were recalculated

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# exp_long
# ---------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("exp_long")
$wsExp.Range("D696").Value = 1676.31878603
$wsExp.Range("D726").Value = 1391.12864775
$wsExp.Range("D733").Value = 1834.46138382

# ---------------------------------------------------------------------
# Table 1 (top expenditure categories, sorted desc by FY22 $ billions)
# ---------------------------------------------------------------------
$wsT1 = $wb.Worksheets.Item("Table 1")

$wsT1.Range("A11").Value = "Il State Toll Highway Auth"
$wsT1.Range("B11").Value = 2.1
$wsT1.Range("C11").Value = 7.06
$wsT1.Range("D11").Value = 7.54

$wsT1.Range("A12").Value = "Debt Service"
$wsT1.Range("B12").Value = 2
$wsT1.Range("C12").Value = -0.83
$wsT1.Range("D12").Value = 6.11

$wsT1.Range("A13").Value = "Revenue"
$wsT1.Range("B13").Value = 1.8
$wsT1.Range("C13").Value = 9.29
$wsT1.Range("D13").Value = 6.33

$wsT1.Range("C16").Value = -17.01
$wsT1.Range("D16").Value = 4.67

# ---------------------------------------------------------------------
# Table 4.b (expenditure category trend table)
# ---------------------------------------------------------------------
$wsT4b = $wb.Worksheets.Item("Table 4.b")

# Row 8 - Commerce And Economic Opportunity
$wsT4b.Range("B8").Value = -17.01
$wsT4b.Range("C8").Value = 49.77
$wsT4b.Range("D8").Value = 34.15
$wsT4b.Range("E8").Value = 16.46
$wsT4b.Range("F8").Value = 3.08
$wsT4b.Range("G8").Value = 4.67

# Row 27 - Revenue
$wsT4b.Range("B27").Value = 9.29
$wsT4b.Range("C27").Value = 27.63
$wsT4b.Range("D27").Value = 45.21
$wsT4b.Range("E27").Value = 30.21
$wsT4b.Range("F27").Value = 13.84
$wsT4b.Range("G27").Value = 6.33

# ---------------------------------------------------------------------
# year_totals - rebuilt table, rows now ordered most-recent-year-first
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("year_totals")

$ws7.Range("A2").Value = 2022
$ws7.Range("B2").Value = 99785.57583777
$ws7.Range("C2").Value = 104536.89862629
$ws7.Range("D2").Value = 4751

$ws7.Range("A3").Value = 2021
$ws7.Range("B3").Value = 92807.10818869
$ws7.Range("C3").Value = 91577.64567625
$ws7.Range("D3").Value = -1229

$ws7.Range("A4").Value = 2020
$ws7.Range("B4").Value = 81574.30708322
$ws7.Range("C4").Value = 74622.7453303
$ws7.Range("D4").Value = -6952

$ws7.Range("A5").Value = 2019
$ws7.Range("B5").Value = 74383.59556887
$ws7.Range("C5").Value = 72152.86792715
$ws7.Range("D5").Value = -2231

$ws7.Range("A6").Value = 2018
$ws7.Range("B6").Value = 74942.56778491
$ws7.Range("C6").Value = 70256.5668322
$ws7.Range("D6").Value = -4686

$ws7.Range("A7").Value = 2017
$ws7.Range("B7").Value = 71704.78677854
$ws7.Range("C7").Value = 60945.18463144
$ws7.Range("D7").Value = -10760

$ws7.Range("A8").Value = 2016
$ws7.Range("B8").Value = 63909.28178688
$ws7.Range("C8").Value = 61806.01279253
$ws7.Range("D8").Value = -2103

$ws7.Range("A9").Value = 2015
$ws7.Range("B9").Value = 69920.57755159
$ws7.Range("C9").Value = 63882.73647204
$ws7.Range("D9").Value = -6038

$ws7.Range("A10").Value = 2014
$ws7.Range("B10").Value = 66941.54371749
$ws7.Range("C10").Value = 62519.59401338
$ws7.Range("D10").Value = -4422

$ws7.Range("A11").Value = 2013
$ws7.Range("B11").Value = 63261.01592636
$ws7.Range("C11").Value = 60502.20379116
$ws7.Range("D11").Value = -2759

$ws7.Range("A12").Value = 2012
$ws7.Range("B12").Value = 59831.151018
$ws7.Range("C12").Value = 55452.47810214
$ws7.Range("D12").Value = -4379

$ws7.Range("A13").Value = 2011
$ws7.Range("B13").Value = 60403.66316511
$ws7.Range("C13").Value = 51719.80617799
$ws7.Range("D13").Value = -8684

$ws7.Range("A14").Value = 2010
$ws7.Range("B14").Value = 59247.71610651
$ws7.Range("C14").Value = 46059.51698249
$ws7.Range("D14").Value = -13188

$ws7.Range("A15").Value = 2009
$ws7.Range("B15").Value = 56721.04766907
$ws7.Range("C15").Value = 47822.53037459
$ws7.Range("D15").Value = -8899

$ws7.Range("A16").Value = 2008
$ws7.Range("B16").Value = 54138.63848686
$ws7.Range("C16").Value = 50213.47771324
$ws7.Range("D16").Value = -3925

$ws7.Range("A17").Value = 2007
$ws7.Range("B17").Value = 51098.59908858
$ws7.Range("C17").Value = 48033.24657002
$ws7.Range("D17").Value = -3065

$ws7.Range("A18").Value = 2006
$ws7.Range("B18").Value = 48028.45089847
$ws7.Range("C18").Value = 44700.58108122
$ws7.Range("D18").Value = -3328

$ws7.Range("A19").Value = 2005
$ws7.Range("B19").Value = 45331.21505246
$ws7.Range("C19").Value = 42865.85989889
$ws7.Range("D19").Value = -2465

$ws7.Range("A20").Value = 2004
$ws7.Range("B20").Value = 52980.20713006
$ws7.Range("C20").Value = 40856.23669512
$ws7.Range("D20").Value = -12124

$ws7.Range("A21").Value = 2003
$ws7.Range("B21").Value = 42567.13604378
$ws7.Range("C21").Value = 36805.69976915
$ws7.Range("D21").Value = -5761

$ws7.Range("A22").Value = 2002
$ws7.Range("B22").Value = 42014.32484476
$ws7.Range("C22").Value = 36825.92689326
$ws7.Range("D22").Value = -5188

$ws7.Range("A23").Value = 2001
$ws7.Range("B23").Value = 40300.24229108
$ws7.Range("C23").Value = 37147.74155936
$ws7.Range("D23").Value = -3153

$ws7.Range("A24").Value = 2000
$ws7.Range("B24").Value = 37283.0484234
$ws7.Range("C24").Value = 35846.01278232
$ws7.Range("D24").Value = -1437

$ws7.Range("A25").Value = 1999
$ws7.Range("B25").Value = 33804.96576153
$ws7.Range("C25").Value = 33030.24759485
$ws7.Range("D25").Value = -775

$ws7.Range("A26").Value = 1998
$ws7.Range("B26").Value = 31218.4556557
$ws7.Range("C26").Value = 31264.6818251
$ws7.Range("D26").Value = 46
